$d = $word.ActiveDocument

# --- Step 1: the actual text change -----------------------------------
# Swap "learning" for "classification" in "This problem is a supervised
# learning. ". Any text mutation re-coalesces every run in the paragraph
# it touches, so this has to happen *before* the bookmark/run-split work
# below, otherwise the splits we create here would just get merged away.
$find = $d.Content
$find.Find.Execute("learning", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wordStart = $find.Start
$wordEnd = $find.End
$wordRng = $d.Range($wordStart, $wordEnd)
$wordRng.Text = "classification"

# --- Step 2: recompute where things live now ---------------------------
$sentence = $d.Content
$sentence.Find.Execute("This problem is a supervised classification. ", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sentenceStart = $sentence.Start
$thisEnd = $sentenceStart + "This ".Length
$classificationEnd = $sentenceStart + "This problem is a supervised classification".Length

# --- Step 3: restore the original run break between the first sentence
#     ("...graduate. ") and "This ..." -----------------------------------
# Dropping a bookmark and deleting it again splits the run at that point
# without re-merging anything and without leaving a stray bookmark
# behind -- the coalescing pass only runs on text mutations, not on
# bookmark add/delete.
$d.Bookmarks.Add("zzTempSplitA", $d.Range($sentenceStart, $sentenceStart)) | Out-Null
$d.Bookmarks("zzTempSplitA").Delete() | Out-Null

# --- Step 4: split "This " away from "problem is a supervised
#     classification" ----------------------------------------------------
$d.Bookmarks.Add("zzTempSplitB", $d.Range($thisEnd, $thisEnd)) | Out-Null
$d.Bookmarks("zzTempSplitB").Delete() | Out-Null

# --- Step 5: drop the "_GoBack" bookmark right after "classification" --
# Word only ever keeps a single "_GoBack" bookmark, so re-adding it here
# both places it at the end of our edit (splitting that run off from the
# trailing ". ") and removes the old one that used to sit at the end of
# the document.
$d.Bookmarks.Add("_GoBack", $d.Range($classificationEnd, $classificationEnd)) | Out-Null
